# Reference latest existing annex instead of original contract in new annexes.
#
# Replaces the hard-coded "Ugovor(a/u) br. {{ broj_ugovora }}" style references
# with generic "referentni_*" placeholders that can point either at the
# original contract or, when present, the latest existing annex.

$d = $word.ActiveDocument

# 1) Title line: "Anex br. {{ broj_aneksa }} Ugovora br. {{ broj_ugovora }}"
$d.Content.Find.Execute(
    "Anex br. {{ broj_aneksa }} Ugovora br. {{ broj_ugovora }}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Anex br. {{ broj_aneksa }} {{ referentni_naziv_gen }} br. {{ referentni_broj }}",
    2)

# 2) Recital: "... dana {{ datum_ugovora }} sklopile Ugovor br. {{ broj_ugovora }} o servisiranju ..."
$d.Content.Find.Execute(
    "Ugovorne strane suglasno utvrđuju da su dana {{ datum_ugovora }} sklopile Ugovor br. {{ broj_ugovora }} o servisiranju i održavanju informacijskog sustava",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ugovorne strane suglasno utvrđuju da su dana {{ datum_referentnog }} sklopile {{ referentni_naziv_nom }} br. {{ referentni_broj }} o servisiranju i održavanju informacijskog sustava",
    2)

# 3) "... odredba čl. 3. gore navedenog Ugovora na način ..."
$d.Content.Find.Execute(
    "Ugovorne strane suglasno utvrđuju da se mijenja odredba čl. 3. gore navedenog Ugovora na način da ista sada glasi:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ugovorne strane suglasno utvrđuju da se mijenja odredba čl. 3. gore navedenog {{ referentni_naziv_gen }} na način da ista sada glasi:",
    2)

# 4) "... odredba čl. 4 st. 4 gore navedenog Ugovora na način ..."
$d.Content.Find.Execute(
    "Ugovorne strane suglasno utvrđuju da se mijenja odredba čl. 4 st. 4 gore navedenog Ugovora na način da ista sada glasi:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ugovorne strane suglasno utvrđuju da se mijenja odredba čl. 4 st. 4 gore navedenog {{ referentni_naziv_gen }} na način da ista sada glasi:",
    2)

# 5) "... odredba Priloga 2. gore navedenog Ugovora na način ..."
$d.Content.Find.Execute(
    "Ugovorne strane suglasno utvrđuju da se mijenja odredba Priloga 2. gore navedenog Ugovora na način da ista sada glasi:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ugovorne strane suglasno utvrđuju da se mijenja odredba Priloga 2. gore navedenog {{ referentni_naziv_gen }} na način da ista sada glasi:",
    2)

# 6) "... ostale odredbe Ugovora {{ broj_ugovora }} o servisiranju i održavanju ... nepromijenjene."
$d.Content.Find.Execute(
    "Ugovorne strane suglasno utvrđuju da ostale odredbe Ugovora {{ broj_ugovora }} o servisiranju i održavanju informacijskog sustava ostaju nepromijenjene.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ugovorne strane suglasno utvrđuju da ostale odredbe {{ referentni_naziv_gen }} {{ referentni_broj }} o servisiranju i održavanju informacijskog sustava ostaju nepromijenjene.",
    2)
